$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("G11").Value = 3.85
$ws.Range("H11").Value = 3.15
$ws.Range("I11").Value = 1.95
$ws.Range("J11").Value = 4.35
$ws.Range("K11").Value = 2.02
$ws.Range("L11").Value = 2.55
$ws.Range("N11").Value = 6.55
$ws.Range("O11").Value = 1.4
$ws.Range("P11").Value = 2.5
$ws.Range("Q11").Value = 2.18
$ws.Range("R11").Value = 1.53
$ws.Range("S11").Value = 1.45
$ws.Range("T11").Value = 2.37
$ws.Range("U11").Value = 1.98
$ws.Range("V11").Value = 1.65
$ws.Range("W11").Value = 8.75
$ws.Range("X11").Value = 19
$ws.Range("Y11").Value = 13.5
$ws.Range("Z11").Value = 60
$ws.Range("AA11").Value = 45
$ws.Range("AB11").Value = 55
$ws.Range("AC11").Value = 7.4
$ws.Range("AD11").Value = 6.2
$ws.Range("AE11").Value = 18
$ws.Range("AF11").Value = 110
$ws.Range("AG11").Value = 5.9
$ws.Range("AH11").Value = 8.25
$ws.Range("AI11").Value = 8.75
$ws.Range("AJ11").Value = 16.5
$ws.Range("AK11").Value = 17.5
$ws.Range("AL11").Value = 35
$ws.Range("AN11").Value = 5.5
$ws.Range("AO11").Value = 22
$ws.Range("AP11").Value = 32
$ws.Range("AQ11").Value = 120
$ws.Range("AR11").Value = 200
$ws.Range("AS11").Value = 500
$ws.Range("AT11").Value = 2.35
$ws.Range("AU11").Value = 7.7
$ws.Range("AW11").Value = 3.65
$ws.Range("AX11").Value = 9.75
$ws.Range("AY11").Value = 21
$ws.Range("AZ11").Value = 37
$ws.Range("BA11").Value = 80
$ws.Range("BB11").Value = 300

# Row 13
$ws.Range("N13").Value = 15

# Row 15
$ws.Range("AM15").Value = 1250

# Row 16
$ws.Range("G16").Value = 2.4
$ws.Range("I16").Value = 2.8
$ws.Range("N16").Value = 8.5
$ws.Range("S16").Value = 1.44
$ws.Range("T16").Value = 2.63
$ws.Range("W16").Value = 7.5
$ws.Range("AC16").Value = 8.5
$ws.Range("AG16").Value = 8
$ws.Range("AH16").Value = 13
$ws.Range("AL16").Value = 34
$ws.Range("AM16").Value = 301
$ws.Range("AN16").Value = 4.5
$ws.Range("AT16").Value = 2.63

# Row 18
$ws.Range("G18").Value = 2.6
$ws.Range("H18").Value = 3.35
$ws.Range("N18").Value = 6.8
$ws.Range("O18").Value = 1.36
$ws.Range("P18").Value = 2.92
$ws.Range("Q18").Value = 2.07
$ws.Range("R18").Value = 1.7
$ws.Range("S18").Value = 1.42
$ws.Range("T18").Value = 2.67
$ws.Range("U18").Value = 1.85
$ws.Range("V18").Value = 1.87
$ws.Range("W18").Value = 7.8
$ws.Range("Z18").Value = 28
$ws.Range("AA18").Value = 23
$ws.Range("AC18").Value = 6.8
$ws.Range("AE18").Value = 15.5
$ws.Range("AF18").Value = 80
$ws.Range("AG18").Value = 7.7
$ws.Range("AH18").Value = 11.75
$ws.Range("AL18").Value = 35
$ws.Range("AM18").Value = 700
$ws.Range("AP18").Value = 23
$ws.Range("AT18").Value = 2.67
$ws.Range("AV18").Value = 70
$ws.Range("AZ18").Value = 60

# Row 19
$ws.Range("M19").Value = 1.11
$ws.Range("N19").Value = 6.5

# Row 21
$ws.Range("Q21").Value = 2.08
$ws.Range("R21").Value = 1.73

# Row 36
$ws.Range("M36").Value = 1.02
$ws.Range("N36").Value = 9.45
